# The source edit swaps ppt/theme/theme1.xml (the slide master's theme,
# previously the custom "Integral" palette) with ppt/theme/theme2.xml (the
# notes master's theme, the stock "Office Theme" palette) -- i.e. after the
# change the deck's slide-master theme uses the standard "Office Theme"
# 12-colour palette instead of "Integral" (font scheme / format scheme are
# byte-identical between the two themes, so only the colour scheme differs).
#
# Reach that colour scheme (12 slots, in fixed order: dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink)
# through any slide and overwrite every slot with the target "Office Theme"
# values. PowerPoint COM reports/accepts colours as 0xBBGGRR (blue, green,
# red) decimal longs via ThemeColor.RGB.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target "Office Theme" palette, expressed as 0xBBGGRR longs, in
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order (indices 1-12).
$officeThemeBgr = @(
    0,          # 1  dk1      RGB 000000
    16777215,   # 2  lt1      RGB FFFFFF
    6968388,    # 3  dk2      RGB 44546A
    15132391,   # 4  lt2      RGB E7E6E6
    13998939,   # 5  accent1  RGB 5B9BD5
    3243501,    # 6  accent2  RGB ED7D31
    10855845,   # 7  accent3  RGB A5A5A5
    49407,      # 8  accent4  RGB FFC000
    12874308,   # 9  accent5  RGB 4472C4
    4697456,    # 10 accent6  RGB 70AD47
    12673797,   # 11 hlink    RGB 0563C1
    7491477     # 12 folHlink RGB 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeThemeBgr[$i - 1]
}
